$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.430.60"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.529.03"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'606.49"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'143.33"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("D7").Value = "3.529.80"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("D11").Value = "'7.73"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "4.119.17"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("E14").Value = "  -6.47%  "
$ws.Range("D15").Value = "'28.76"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "3.521.17"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "66.221.74"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'10.87"
$ws.Range("E19").Value = "  -5.47%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "'14.67"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'424.63"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'0.593"
$ws.Range("E23").Value = "  -4.18%  "
$ws.Range("D24").Value = "'77.20"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").Value = "3.659.32"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'0.0000114"
$ws.Range("E27").Value = "  -6.52%  "
$ws.Range("D28").Value = "'7.97"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "'2.46"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.156"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").Value = "3.525.60"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "'24.34"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -8.83%  "
$ws.Range("D37").Value = "'7.56"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("D39").Value = "'174.59"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "'5.24"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("E41").Value = "  -3.93%  "
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("D43").Value = "'0.860"
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("D44").Value = "'45.35"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  -7.39%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "'2.37"
$ws.Range("E47").Value = "  -7.51%  "
$ws.Range("D48").Value = "'7.10"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'23.03"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("E50").Value = "  -5.62%  "
$ws.Range("D51").Value = "'0.911"
$ws.Range("E51").Value = "  -4.38%  "
